# Adds the "Facebook Multilingual Task Oriented Dataset" row to the survey sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19
$prevRow = 18

# Copy the formatting of the previous data row onto the new row first, so the
# new cells share the existing style (s="1") rather than generating new ones.
$ws.Range("A$prevRow`:I$prevRow").Copy()
$ws.Range("A$row`:I$row").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value  = "Facebook Multilingual Task Oriented Dataset"
$ws.Cells.Item($row, 2).Value  = "1. (Faceboook)  We release a dataset of around 57k annotated utterances`nin English (43k), Spanish (8.6k) and Thai (5k) for three task oriented domains … ALARM,`nREMINDER, and WEATHER.`n2. For cross-lingual natural language understanding"
$ws.Cells.Item($row, 3).Value  = "Download: https://fb.me/multilingual_task_oriented_data `nPaper: https://arxiv.org/pdf/1810.13327.pdf "
$ws.Cells.Item($row, 4).Value  = "S"
$ws.Cells.Item($row, 5).Value  = "Task Oriented"
$ws.Cells.Item($row, 6).Value  = "3 Domains: Alarm, Reminder, Weather`n3 Languages: English, Spanish, Thai"
$ws.Cells.Item($row, 7).Value  = "Yes"
$ws.Cells.Item($row, 8).Value  = "English Train: 30,521`nEnglish Dev: 4,181`nEnglish Test: 8,621`nSpanish Train: 3,617`nSpanish Dev: 1,983`nSpanish Test: 3,043`nThai Train: 2,156`nThai Dev: 1,235`nThai Test:  1,692"
$ws.Cells.Item($row, 9).Value  = "Slot`nIntent"

# row height to fit the new (taller) content, as in the saved workbook
$ws.Rows.Item($row).RowHeight = 169

# update the active selection to the newly added row, as in the saved workbook
[void]$ws.Range("A$row").Select()
